# Update header labels from cm to in
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "box_length_in"
$ws.Range("F1").Value = "box_width_in"
$ws.Range("G1").Value = "box_height_in"

# Convert box dimension values (cm -> in) for each data row
$ws.Range("E2").Value = 15.7
$ws.Range("F2").Value = 11.8
$ws.Range("G2").Value = 9.8

$ws.Range("E3").Value = 15.7
$ws.Range("F3").Value = 11.8
$ws.Range("G3").Value = 9.8

$ws.Range("E4").Value = 17.7
$ws.Range("F4").Value = 13.8
$ws.Range("G4").Value = 11

$ws.Range("E5").Value = 17.7
$ws.Range("F5").Value = 13.8
$ws.Range("G5").Value = 11

$ws.Range("E6").Value = 19.7
$ws.Range("F6").Value = 13.8
$ws.Range("G6").Value = 7.9

$ws.Range("E7").Value = 18.9
$ws.Range("F7").Value = 12.6
$ws.Range("G7").Value = 8.7

$ws.Range("E8").Value = 19.7
$ws.Range("F8").Value = 13.8
$ws.Range("G8").Value = 9.8

$ws.Range("E9").Value = 16.5
$ws.Range("F9").Value = 11.8
$ws.Range("G9").Value = 7.9
